# Generate Report for handback
# Refresh the "Latest Handoff Datetime" (column D) and "Latest Handback
# DateTime" (column G) for the most recently processed source file (row 2)
# on both the "zh-cn" and "de-de" status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-09 03:55:08"
$wsZhCn.Range("G2").Value = "2016-01-09 03:55:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-09 03:55:17"
$wsDeDe.Range("G2").Value = "2016-01-09 03:56:07"
